# Updated cryptos list on Sat Jun  3 23:57:12 UTC 2023 with GitHub Actions
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "27.127.47"
$ws.Range("E2").Value = "  -0.45%  "
$ws.Range("D3").Value = "1.896.10"
$ws.Range("E3").Value = "  -0.60%  "
$ws.Range("E4").Value = "  +0.28%  "
$ws.Range("D5").Value = "'307.50"
$ws.Range("E5").Value = "  +0.01%  "
$ws.Range("E6").Value = "  +0.22%  "
$ws.Range("D7").Value = "'0.5201"
$ws.Range("E7").Value = "  -1.03%  "
$ws.Range("D8").Value = "'0.3770"
$ws.Range("E8").Value = "  -0.37%  "
$ws.Range("D9").Value = "'0.07288"
$ws.Range("E9").Value = "  +0.38%  "
$ws.Range("E10").Value = "  -0.31%  "
$ws.Range("D11").Value = "'0.9017"
$ws.Range("E11").Value = "  +0.16%  "
$ws.Range("D12").Value = "'0.08204"
$ws.Range("E12").Value = "  -0.78%  "
$ws.Range("D13").Value = "1.943.52"
$ws.Range("E13").Value = "  +1.50%  "
$ws.Range("D14").Value = "'96.14"
$ws.Range("E14").Value = "  +0.84%  "
$ws.Range("D15").Value = "'5.339"
$ws.Range("E15").Value = "  +1.12%  "
$ws.Range("D16").Value = "'1.004"
$ws.Range("D17").Value = "'0.000008622"
$ws.Range("E17").Value = "  +0.19%  "
$ws.Range("E18").Value = "  +0.88%  "
$ws.Range("D19").Value = "'1.002"
$ws.Range("E19").Value = "  +0.25%  "
$ws.Range("D20").Value = "27.169.52"
$ws.Range("E20").Value = "  -0.42%  "
$ws.Range("D21").Value = "'5.098"
$ws.Range("E21").Value = "  +0.66%  "
$ws.Range("D23").Value = "'6.433"
$ws.Range("E23").Value = "  -0.36%  "
$ws.Range("D24").Value = "'148.98"
$ws.Range("E24").Value = "  +2.06%  "
$ws.Range("D25").Value = "'2.313"
$ws.Range("E25").Value = "  +0.57%  "
$ws.Range("B26").Value = "Toncoin"
$ws.Range("C26").Value = "https://coinranking.com/coin/67YlI0K1b+toncoin-ton"
$ws.Range("D26").Value = "'1.744"
$ws.Range("E26").Value = "  +0.00%  "
$ws.Range("B27").Value = "EthereumClassic"
$ws.Range("C27").Value = "https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc"
$ws.Range("D27").Value = "'18.18"
$ws.Range("E27").Value = "  +0.12%  "
$ws.Range("D29").Value = "'4.808"
$ws.Range("E29").Value = "  -0.07%  "
$ws.Range("E30").Value = "  -2.05%  "
$ws.Range("E31").Value = "  +0.15%  "
$ws.Range("D32").Value = "'0.7964"
$ws.Range("E32").Value = "  -0.82%  "
$ws.Range("D33").Value = "'0.05033"
$ws.Range("E33").Value = "  -1.05%  "
$ws.Range("D34").Value = "'1.219"
$ws.Range("E34").Value = "  -1.60%  "
$ws.Range("E35").Value = "  +1.89%  "
$ws.Range("E36").Value = "  +0.23%  "
$ws.Range("D37").Value = "'2.611"
$ws.Range("E37").Value = "  +1.64%  "
$ws.Range("D38").Value = "'0.5737"
$ws.Range("E38").Value = "  -0.18%  "
$ws.Range("E39").Value = "  +1.11%  "
$ws.Range("D40").Value = "'1.079"
$ws.Range("E40").Value = "  +0.41%  "
$ws.Range("D41").Value = "'9.000"
$ws.Range("E41").Value = "  -0.58%  "
$ws.Range("D42").Value = "'6.558"
$ws.Range("E42").Value = "  -1.19%  "
$ws.Range("D43").Value = "'116.51"
$ws.Range("E43").Value = "  -1.72%  "
$ws.Range("E44").Value = "  -0.34%  "
$ws.Range("D45").Value = "'0.4875"
$ws.Range("E45").Value = "  +0.71%  "
$ws.Range("E46").Value = "  +0.18%  "
$ws.Range("D47").Value = "'10.13"
$ws.Range("E47").Value = "  -0.49%  "
$ws.Range("D48").Value = "'1.620"
$ws.Range("E48").Value = "  +0.41%  "
$ws.Range("D49").Value = "'38.39"
$ws.Range("E49").Value = "  +2.01%  "
$ws.Range("D50").Value = "'63.81"
$ws.Range("E50").Value = "  +0.19%  "
$ws.Range("D51").Value = "'0.05930"
$ws.Range("E51").Value = "  -0.35%  "
